# Backward-elimination OLS summary blobs were re-generated on a later run;
# only the embedded "Date:" / "Time:" stamps inside the big summary text
# blocks (column B, row 2 of each step sheet) actually changed.
#
# Sheets 1-19 (tab "41" .. tab "23") were (re)written at 23:19:09,
# sheets 20-28 (tab "22" .. tab "14") a moment later at 23:19:10.

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 29 Dec 2019"
$newDate = "Wed, 01 Jan 2020"
$oldTime = "16:11:33"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    if ($i -le 19) {
        $newTime = "23:19:09"
    } else {
        $newTime = "23:19:10"
    }

    $ws.Cells.Replace($oldDate, $newDate)
    $ws.Cells.Replace($oldTime, $newTime)
}
